$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "K" column (column G) values for rows 2-5 to reflect the
# regenerated save data (K replaces the old Strike# values).
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
